$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Oppvisningsstevne iskanten 9 februar 2026"
$ws.Range("D5").Value = "17:07:45"
$ws.Range("E5").Value = "Aksel Eriksen"
$ws.Range("F5").Value = "Loddefjord IL"
$ws.Range("B6").Value = "17:07:45"
$ws.Range("D6").Value = "17:11:30"
$ws.Range("E6").Value = "Frida Pasko Hansen"
$ws.Range("F6").Value = "Loddefjord IL"
$ws.Range("B7").Value = "17:11:30"
$ws.Range("D7").Value = "17:15:15"
$ws.Range("E7").Value = "Angela Chen"
$ws.Range("F7").Value = "Fana Idrettslag"
$ws.Range("B8").Value = "17:15:15"
$ws.Range("D8").Value = "17:19:00"
$ws.Range("E8").Value = "Mille Isabell Steen Rein"
$ws.Range("F8").Value = "Loddefjord IL"
$ws.Range("B9").Value = "17:19:00"
$ws.Range("D9").Value = "17:22:45"
$ws.Range("E9").Value = "Aylin Morseth"
$ws.Range("F9").Value = "Fana Idrettslag"
$ws.Range("B10").Value = "17:22:45"
$ws.Range("D10").Value = "17:26:30"
$ws.Range("E10").Value = "Camilla Tveit"
$ws.Range("F10").Value = "Loddefjord IL"
$ws.Range("B11").Value = "17:26:30"
$ws.Range("D11").Value = "17:30:15"
$ws.Range("E11").Value = "Sara Barbro Kyte"
$ws.Range("F11").Value = "Bergen Kunstløpklubb"
$ws.Range("B12").Value = "17:30:15"
$ws.Range("D12").Value = "17:34:00"
$ws.Range("E12").Value = "Elena Sophia Sandnes-Strømmen"
$ws.Range("B13").Value = "ca. 17:34:00"
$ws.Range("D13").Value = "17:38:00"
$ws.Range("B14").Value = "17:38:00"
$ws.Range("D14").Value = "17:41:45"
$ws.Range("E14").Value = "Frida Lovisa Østerberg"
$ws.Range("B15").Value = "17:41:45"
$ws.Range("D15").Value = "17:45:30"
$ws.Range("E15").Value = "Amanda Ansnes Lima"
$ws.Range("B16").Value = "17:45:30"
$ws.Range("D16").Value = "17:49:15"
$ws.Range("E16").Value = "Aurelia Landschulze"
$ws.Range("F16").Value = "Fana Idrettslag"
$ws.Range("B17").Value = "17:49:15"
$ws.Range("D17").Value = "17:53:00"
$ws.Range("E17").Value = "Valentina Pinker-Spilde"
$ws.Range("F17").Value = "Fana Idrettslag"
$ws.Range("B18").Value = "17:53:00"
$ws.Range("D18").Value = "17:56:45"
$ws.Range("E18").Value = "Eira Olava Bortne Ludvigsen"
$ws.Range("B19").Value = "17:56:45"
$ws.Range("D19").Value = "18:00:30"
$ws.Range("E19").Value = "Sarolt Szofia Papdi"
$ws.Range("B20").Value = "18:00:30"
$ws.Range("D20").Value = "18:04:15"
$ws.Range("E20").Value = "Mie Mariell Sævereid"
$ws.Range("F20").Value = "Bergen Kunstløpklubb"
$ws.Range("B21").Value = "18:04:15"
$ws.Range("D21").Value = "18:08:00"
$ws.Range("E21").Value = "Hennie Markestad"
$ws.Range("F21").Value = "Bergen Kunstløpklubb"
$ws.Range("B22").Value = "ca. 18:08:00"
$ws.Range("D22").Value = "18:12:00"
$ws.Range("B23").Value = "18:12:00"
$ws.Range("D23").Value = "18:15:45"
$ws.Range("E23").Value = "Eleanora Egle"
$ws.Range("F23").Value = "Loddefjord IL"
$ws.Range("B24").Value = "18:15:45"
$ws.Range("D24").Value = "18:19:30"
$ws.Range("E24").Value = "Patricija Levickaite"
$ws.Range("B25").Value = "18:19:30"
$ws.Range("D25").Value = "18:23:15"
$ws.Range("E25").Value = "Yuewei Li"
$ws.Range("F25").Value = "Bergen Kunstløpklubb"
$ws.Range("B26").Value = "18:23:15"
$ws.Range("D26").Value = "18:27:00"
$ws.Range("E26").Value = "Hanna Wangsuk Tveita"
$ws.Range("B27").Value = "18:27:00"
$ws.Range("D27").Value = "18:30:45"
$ws.Range("E27").Value = "Leah Kalvik"
$ws.Range("B28").Value = "18:30:45"
$ws.Range("D28").Value = "18:34:30"
$ws.Range("E28").Value = "Emilie Morseth"
$ws.Range("F28").Value = "Fana Idrettslag"
$ws.Range("B29").Value = "18:34:30"
$ws.Range("D29").Value = "18:38:15"
$ws.Range("E29").Value = "Frida Qianlu He"
$ws.Range("F29").Value = "Loddefjord IL"
$ws.Range("B30").Value = "18:38:15"
$ws.Range("D30").Value = "18:42:00"
$ws.Range("E30").Value = "Anne Kristoffersen"
$ws.Range("A32").Value = "Generert 01.02.2026 08:45 • OLES_LAPTOP"
